$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---- Column G: 12 changed rows ----
$ws.Range("G68").Value = "saveModalDialogTextByLocator(var,locator)"
$ws.Range("G85").Value = "typeKeys(os,keystrokes)"
$ws.Range("G86").Value = "typeTextArea(name,text1,text2,text3,text4)"
$ws.Range("G87").Value = "typeTextBox(name,text1,text2,text3,text4)"
$ws.Range("G88").Value = "useApp(appId)"
$ws.Range("G89").Value = "useForm(formName)"
$ws.Range("G90").Value = "useHierTable(var,name)"
$ws.Range("G91").Value = "useList(var,name)"
$ws.Range("G92").Value = "useTable(var,name)"
$ws.Range("G93").Value = "useTableRow(var,row)"
$ws.Range("G94").Value = "waitFor(name,maxWaitMs)"
$ws.Range("G95").Value = "waitForLocator(locator,maxWaitMs)"
# ---- Column J: 1 changed rows ----
$ws.Range("J7").Value = "saveDiff(var,baseline,actual)"
# ---- Column M: 6 changed rows ----
$ws.Range("M12").Value = "compact(var,json,removeEmpty)"
$ws.Range("M13").Value = "fromCsv(csv,header,jsonFile)"
$ws.Range("M14").Value = "minify(json,var)"
$ws.Range("M15").Value = "storeCount(json,jsonpath,var)"
$ws.Range("M16").Value = "storeValue(json,jsonpath,var)"
$ws.Range("M17").Value = "storeValues(json,jsonpath,var)"
# ---- Column Y: 80 changed rows ----
$ws.Range("Y50").Value = "clickAll(locator)"
$ws.Range("Y51").Value = "clickAndWait(locator,waitMs)"
$ws.Range("Y52").Value = "clickByLabel(label)"
$ws.Range("Y53").Value = "clickByLabelAndWait(label,waitMs)"
$ws.Range("Y54").Value = "clickOffset(locator,x,y)"
$ws.Range("Y55").Value = "clickWithKeys(locator,keys)"
$ws.Range("Y56").Value = "close()"
$ws.Range("Y57").Value = "closeAll()"
$ws.Range("Y58").Value = "deselect(locator,text)"
$ws.Range("Y59").Value = "deselectMulti(locator,array)"
$ws.Range("Y60").Value = "dismissInvalidCert()"
$ws.Range("Y61").Value = "dismissInvalidCertPopup()"
$ws.Range("Y62").Value = "doubleClick(locator)"
$ws.Range("Y63").Value = "doubleClickAndWait(locator,waitMs)"
$ws.Range("Y64").Value = "doubleClickByLabel(label)"
$ws.Range("Y65").Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Range("Y66").Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Range("Y67").Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Range("Y68").Value = "editLocalStorage(key,value)"
$ws.Range("Y69").Value = "executeScript(var,script)"
$ws.Range("Y70").Value = "focus(locator)"
$ws.Range("Y71").Value = "goBack()"
$ws.Range("Y72").Value = "goBackAndWait()"
$ws.Range("Y73").Value = "maximizeWindow()"
$ws.Range("Y74").Value = "mouseOver(locator)"
$ws.Range("Y75").Value = "open(url)"
$ws.Range("Y76").Value = "openAndWait(url,waitMs)"
$ws.Range("Y77").Value = "openHttpBasic(url,username,password)"
$ws.Range("Y78").Value = "openIgnoreTimeout(url)"
$ws.Range("Y79").Value = "refresh()"
$ws.Range("Y80").Value = "refreshAndWait()"
$ws.Range("Y81").Value = "resizeWindow(width,height)"
$ws.Range("Y82").Value = "rightClick(locator)"
$ws.Range("Y83").Value = "saveAllWindowIds(var)"
$ws.Range("Y84").Value = "saveAllWindowNames(var)"
$ws.Range("Y85").Value = "saveAttribute(var,locator,attrName)"
$ws.Range("Y86").Value = "saveAttributeList(var,locator,attrName)"
$ws.Range("Y87").Value = "saveCount(var,locator)"
$ws.Range("Y88").Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Range("Y89").Value = "saveElement(var,locator)"
$ws.Range("Y90").Value = "saveElements(var,locator)"
$ws.Range("Y91").Value = "saveLocalStorage(var,key)"
$ws.Range("Y92").Value = "saveLocation(var)"
$ws.Range("Y93").Value = "savePageAs(var,sessionIdName,url)"
$ws.Range("Y94").Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Range("Y95").Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Range("Y96").Value = "saveText(var,locator)"
$ws.Range("Y97").Value = "saveTextArray(var,locator)"
$ws.Range("Y98").Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Range("Y99").Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Range("Y100").Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Range("Y101").Value = "saveValue(var,locator)"
$ws.Range("Y102").Value = "saveValues(var,locator)"
$ws.Range("Y103").Value = "scrollElement(locator,xOffset,yOffset)"
$ws.Range("Y104").Value = "scrollLeft(locator,pixel)"
$ws.Range("Y105").Value = "scrollPage(xOffset,yOffset)"
$ws.Range("Y106").Value = "scrollRight(locator,pixel)"
$ws.Range("Y107").Value = "scrollTo(locator)"
$ws.Range("Y108").Value = "select(locator,text)"
$ws.Range("Y109").Value = "selectFrame(locator)"
$ws.Range("Y110").Value = "selectMulti(locator,array)"
$ws.Range("Y111").Value = "selectMultiOptions(locator)"
$ws.Range("Y112").Value = "selectText(locator)"
$ws.Range("Y113").Value = "selectWindow(winId)"
$ws.Range("Y114").Value = "selectWindowAndWait(winId,waitMs)"
$ws.Range("Y115").Value = "selectWindowByIndex(index)"
$ws.Range("Y116").Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Range("Y117").Value = "toggleSelections(locator)"
$ws.Range("Y118").Value = "type(locator,value)"
$ws.Range("Y119").Value = "typeKeys(locator,value)"
$ws.Range("Y120").Value = "uncheckAll(locator)"
$ws.Range("Y121").Value = "unselectAllText()"
$ws.Range("Y122").Value = "upload(fieldLocator,file)"
$ws.Range("Y123").Value = "verifyContainText(locator,text)"
$ws.Range("Y124").Value = "verifyText(locator,text)"
$ws.Range("Y125").Value = "wait(waitMs)"
$ws.Range("Y126").Value = "waitForElementPresent(locator)"
$ws.Range("Y127").Value = "waitForPopUp(winId,waitMs)"
$ws.Range("Y128").Value = "waitForTextPresent(text)"
$ws.Range("Y129").Value = "waitForTitle(text)"
# ---- Column AD: 21 changed rows ----
$ws.Range("AD7").Value = "assertSoap(wsdl,xml)"
$ws.Range("AD8").Value = "assertSoapFaultCode(expected,xml)"
$ws.Range("AD9").Value = "assertSoapFaultString(expected,xml)"
$ws.Range("AD10").Value = "assertValue(xml,xpath,expected)"
$ws.Range("AD11").Value = "assertValues(xml,xpath,array,exactOrder)"
$ws.Range("AD12").Value = "assertWellformed(xml)"
$ws.Range("AD13").Value = "beautify(xml,var)"
$ws.Range("AD14").Value = "clear(xml,xpath,var)"
$ws.Range("AD15").Value = "delete(xml,xpath,var)"
$ws.Range("AD16").Value = "insertAfter(xml,xpath,content,var)"
$ws.Range("AD17").Value = "insertBefore(xml,xpath,content,var)"
$ws.Range("AD18").Value = "minify(xml,var)"
$ws.Range("AD19").Value = "prepend(xml,xpath,content,var)"
$ws.Range("AD20").Value = "replace(xml,xpath,content,var)"
$ws.Range("AD21").Value = "replaceIn(xml,xpath,content,var)"
$ws.Range("AD22").Value = "storeCount(xml,xpath,var)"
$ws.Range("AD23").Value = "storeSoapFaultCode(var,xml)"
$ws.Range("AD24").Value = "storeSoapFaultDetail(var,xml)"
$ws.Range("AD25").Value = "storeSoapFaultString(var,xml)"
$ws.Range("AD26").Value = "storeValue(xml,xpath,var)"
$ws.Range("AD27").Value = "storeValues(xml,xpath,var)"

# ---- Update defined names (ranges) to reflect new extents ----
$wb.Names.Item("desktop").RefersTo = "='#system'!`$G`$2:`$G`$95"
$wb.Names.Item("image").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$17"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

Write-Output "done"
